# "Modification of infectiousness by age"
#
# Inserts three new parameter rows at the top of the `constants` sheet
# (rows 4-6) for the new child-infectiousness multipliers, pushing every
# row from the former row 4 onward down by three. Also restores the
# selection/page-setup cosmetics that Excel re-wrote on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Insert 3 blank rows above the old row 4 ("tb_prop_early_progression").
# Excel clones row 4's formatting (styles 37/38/39) onto the new rows and
# shifts every downstream row reference (dimension, dataValidations, named
# formulas, etc.) automatically.
$ws.Rows("4:6").Insert()

# New parameter: overall child infectiousness multiplier.
$ws.Range("A4").Value = "tb_multiplier_child_infectiousness"
$ws.Range("B4").Value = 1

# New parameter: child infectiousness multiplier, age 0-10.
$ws.Range("A5").Value = "tb_multiplier_child_infectiousness_age0to10"
$ws.Range("B5").Value = 0.1

# New parameter: child infectiousness multiplier, age 10+.
$ws.Range("A6").Value = "tb_multiplier_child_infectiousness_age10up"
$ws.Range("B6").Value = 1

# Restore the sheet's printed page setup (paper size / orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Reset the view: scroll back to the top and select A7 (just below the
# newly-inserted rows), matching the saved selection state.
$ws.Range("A1").Select()
$ws.Range("A7").Select()
